# Updated legacy GSC export data:
# The oldest date row ("2025-11-19", the first data row) is removed from the
# "Chart" sheet. All subsequent rows shift up by one, so what used to be
# "2025-11-20" becomes the new first data row, etc. This mirrors deleting
# worksheet row 2 and letting Excel re-flow the remaining rows upward.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Delete the entire second row (first data row, date 2025-11-19).
# This shifts every row below it up by one, shrinking the used range
# from A1:D89 down to A1:D88.
$ws.Rows.Item(2).Delete()
